$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Espinaca (Vega Monumental Concepción) needs to
# be inserted as row 45, pushing the existing rows 45-61 down to 46-62.
$ws.Rows("45").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 44636
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = 100112012
$ws.Cells.Item(45, 7).Value = "Espinaca"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 220
$ws.Cells.Item(45, 11).Value = 8000
$ws.Cells.Item(45, 12).Value = 9000
$ws.Cells.Item(45, 13).Value = 8545
$ws.Cells.Item(45, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(45, 15).Value = "Región Metropolitana"
$ws.Cells.Item(45, 16).Value = 854
$ws.Cells.Item(45, 17).Value = 10
$ws.Cells.Item(45, 18).Value = "Hortaliza"
